$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "finalidad" in column L, matching the formatting of K1.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("L1").Value = "finalidad"

# Move the active selection to K8 (matches recorded selection in the saved file)
$ws.Range("K8").Select() | Out-Null

$wb.Save()
